$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3-6 (Turno2..Turno5), keeping only the Turno1 data row
$ws.Range("A3:D6").EntireRow.Delete() | Out-Null

# Update the "decorre" value for Turno1 to reflect the new 0/1 availability
$ws.Range("B2").Value = "[1,0,0,0,0,0,0,0,0,0]"

# Update selection to match the saved state
$ws.Range("B3").Select() | Out-Null
